$wb = $excel.ActiveWorkbook

function Set-Cell($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
Set-Cell $ws 2 8 2155   # H2
Set-Cell $ws 2 9 792.5   # I2
Set-Cell $ws 2 10 2700   # J2
Set-Cell $ws 2 11 792.5   # K2
Set-Cell $ws 2 12 2700   # L2
Set-Cell $ws 2 13 -679.5   # M2
Set-Cell $ws 2 14 -2926   # N2
Set-Cell $ws 4 8 341.66666   # H4
Set-Cell $ws 4 9 350   # I4
Set-Cell $ws 4 10 325   # J4
Set-Cell $ws 4 11 350   # K4
Set-Cell $ws 4 12 325   # L4
Set-Cell $ws 4 13 -236   # M4
Set-Cell $ws 4 14 -553   # N4
Set-Cell $ws 5 8 404.6   # H5
Set-Cell $ws 5 9 404.6   # I5
Set-Cell $ws 5 10 0   # J5
Set-Cell $ws 5 11 404.6   # K5
Set-Cell $ws 5 12 0   # L5
Set-Cell $ws 5 13 -289.6   # M5
Set-Cell $ws 6 8 864.55554   # H6
Set-Cell $ws 6 9 476.8   # I6
Set-Cell $ws 6 10 1349.25   # J6
Set-Cell $ws 6 11 1430.4   # K6
Set-Cell $ws 6 12 4047.75   # L6
Set-Cell $ws 6 13 -1318.4   # M6
Set-Cell $ws 6 14 -4271.75   # N6
Set-Cell $ws 8 8 283.92062   # H8
Set-Cell $ws 8 9 180.25   # I8
Set-Cell $ws 8 10 299   # J8
Set-Cell $ws 8 11 540.75   # K8
Set-Cell $ws 8 12 897   # L8
Set-Cell $ws 8 13 -401.75   # M8
Set-Cell $ws 8 14 -1175   # N8
Set-Cell $ws 9 8 883.3333   # H9
Set-Cell $ws 9 9 980   # I9
Set-Cell $ws 9 10 400   # J9
Set-Cell $ws 9 11 980   # K9
Set-Cell $ws 9 12 400   # L9
Set-Cell $ws 9 13 -811   # M9
Set-Cell $ws 9 14 -738   # N9
Set-Cell $ws 38 8 1022.2   # H38
Set-Cell $ws 38 9 152.75   # I38
Set-Cell $ws 38 10 4500   # J38
Set-Cell $ws 38 11 458.25   # K38
Set-Cell $ws 38 12 13500   # L38
Set-Cell $ws 38 13 -86.25   # M38
Set-Cell $ws 38 14 -14244   # N38
Set-Cell $ws 80 8 21380.5   # H80
Set-Cell $ws 80 9 10159.6   # I80
Set-Cell $ws 80 10 32601.4   # J80
Set-Cell $ws 80 11 30478.8   # K80
Set-Cell $ws 80 12 97804.20000000001   # L80
Set-Cell $ws 80 13 -29480.8   # M80
Set-Cell $ws 80 14 -99800.20000000001   # N80
Set-Cell $ws 83 8 21380.5   # H83
Set-Cell $ws 83 9 10159.6   # I83
Set-Cell $ws 83 10 32601.4   # J83
Set-Cell $ws 83 11 91436.40000000001   # K83
Set-Cell $ws 83 12 293412.6   # L83
Set-Cell $ws 83 13 -86444.40000000001   # M83
Set-Cell $ws 83 14 -303396.6   # N83
Set-Cell $ws 100 8 48964.152   # H100
Set-Cell $ws 100 9 52919.5   # I100
Set-Cell $ws 100 10 1500   # J100
Set-Cell $ws 100 11 52919.5   # K100
Set-Cell $ws 100 12 1500   # L100
Set-Cell $ws 100 13 -52378.5   # M100
Set-Cell $ws 100 14 -2582   # N100
Set-Cell $ws 123 8 69998.5   # H123
Set-Cell $ws 123 9 0   # I123
Set-Cell $ws 123 10 69998.5   # J123
Set-Cell $ws 123 11 0   # K123
Set-Cell $ws 123 12 69998.5   # L123
Set-Cell $ws 123 14 -79798.5   # N123
Set-Cell $ws 128 8 115000   # H128
Set-Cell $ws 128 9 0   # I128
Set-Cell $ws 128 10 115000   # J128
Set-Cell $ws 128 11 0   # K128
Set-Cell $ws 128 12 115000   # L128
Set-Cell $ws 128 14 -124960   # N128
Set-Cell $ws 129 8 2118.75   # H129
Set-Cell $ws 129 9 983.3333   # I129
Set-Cell $ws 129 10 2800   # J129
Set-Cell $ws 129 11 2949.9999   # K129
Set-Cell $ws 129 12 8400   # L129
Set-Cell $ws 129 13 2050.0001   # M129
Set-Cell $ws 129 14 -18400   # N129
Set-Cell $ws 132 8 2660805   # H132
Set-Cell $ws 132 9 2993079.5   # I132
Set-Cell $ws 132 10 2610   # J132
Set-Cell $ws 132 11 8979238.5   # K132
Set-Cell $ws 132 12 7830   # L132
Set-Cell $ws 132 13 -8976708.5   # M132
Set-Cell $ws 132 14 -12890   # N132
Set-Cell $ws 135 8 12277.885   # H135
Set-Cell $ws 135 9 914.439   # I135
Set-Cell $ws 135 10 54632.547   # J135
Set-Cell $ws 135 11 8229.950999999999   # K135
Set-Cell $ws 135 12 491692.923   # L135
Set-Cell $ws 135 13 -5694.950999999999   # M135
Set-Cell $ws 135 14 -496762.923   # N135
Set-Cell $ws 137 8 18339.316   # H137
Set-Cell $ws 137 9 29463.273   # I137
Set-Cell $ws 137 10 3043.875   # J137
Set-Cell $ws 137 11 88389.819   # K137
Set-Cell $ws 137 12 9131.625   # L137
Set-Cell $ws 137 13 -85839.819   # M137
Set-Cell $ws 137 14 -14231.625   # N137
Set-Cell $ws 138 8 1518.4833   # H138
Set-Cell $ws 138 9 1079.7255   # I138
Set-Cell $ws 138 10 4004.7778   # J138
Set-Cell $ws 138 11 3239.1765   # K138
Set-Cell $ws 138 12 12014.3334   # L138
Set-Cell $ws 138 13 1900.8235   # M138
Set-Cell $ws 138 14 -22294.3334   # N138

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
Set-Cell $ws 32 8 20605.21   # H32
Set-Cell $ws 32 9 15511.949   # I32
Set-Cell $ws 32 10 87666.5   # J32
Set-Cell $ws 32 11 15511.949   # K32
Set-Cell $ws 32 12 87666.5   # L32
Set-Cell $ws 32 13 -15224.949   # M32
Set-Cell $ws 32 14 -88240.5   # N32
Set-Cell $ws 97 8 1202.9269   # H97
Set-Cell $ws 97 9 772.3913   # I97
Set-Cell $ws 97 10 1753.0555   # J97
Set-Cell $ws 97 11 772.3913   # K97
Set-Cell $ws 97 12 1753.0555   # L97
Set-Cell $ws 97 13 -276.3913   # M97
Set-Cell $ws 97 14 -2745.0555   # N97
Set-Cell $ws 132 8 2133.6135   # H132
Set-Cell $ws 132 9 1909.6552   # I132
Set-Cell $ws 132 10 2566.6   # J132
Set-Cell $ws 132 11 5728.9656   # K132
Set-Cell $ws 132 12 7699.799999999999   # L132
Set-Cell $ws 132 13 -3198.9656   # M132
Set-Cell $ws 132 14 -12759.8   # N132

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
Set-Cell $ws 86 8 1854.8182   # H86
Set-Cell $ws 86 9 2284.3333   # I86
Set-Cell $ws 86 10 1339.4   # J86
Set-Cell $ws 86 11 2284.3333   # K86
Set-Cell $ws 86 12 1339.4   # L86
Set-Cell $ws 86 13 -1161.3333   # M86
Set-Cell $ws 86 14 -3585.4   # N86
Set-Cell $ws 89 8 1854.8182   # H89
Set-Cell $ws 89 9 2284.3333   # I89
Set-Cell $ws 89 10 1339.4   # J89
Set-Cell $ws 89 11 11421.6665   # K89
Set-Cell $ws 89 12 6697   # L89
Set-Cell $ws 89 13 -5805.666499999999   # M89
Set-Cell $ws 89 14 -17929   # N89
Set-Cell $ws 134 8 2679.0144   # H134
Set-Cell $ws 134 9 1368.7637   # I134
Set-Cell $ws 134 10 7483.2666   # J134
Set-Cell $ws 134 11 4106.2911   # K134
Set-Cell $ws 134 12 22449.7998   # L134
Set-Cell $ws 134 13 -1571.2911   # M134
Set-Cell $ws 134 14 -27519.7998   # N134

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
Set-Cell $ws 58 8 1325.9166   # H58
Set-Cell $ws 58 9 1341.4736   # I58
Set-Cell $ws 58 10 1266.8   # J58
Set-Cell $ws 58 11 1341.4736   # K58
Set-Cell $ws 58 12 1266.8   # L58
Set-Cell $ws 58 13 -1138.4736   # M58
Set-Cell $ws 58 14 -1672.8   # N58
Set-Cell $ws 132 8 56746.223   # H132
Set-Cell $ws 132 9 67562.13   # I132
Set-Cell $ws 132 10 2666.6667   # J132
Set-Cell $ws 132 11 202686.39   # K132
Set-Cell $ws 132 12 8000.000100000001   # L132
Set-Cell $ws 132 13 -200156.39   # M132
Set-Cell $ws 132 14 -13060.0001   # N132
Set-Cell $ws 134 8 2136.0571   # H134
Set-Cell $ws 134 9 2096.8262   # I134
Set-Cell $ws 134 10 2211.25   # J134
Set-Cell $ws 134 11 6290.4786   # K134
Set-Cell $ws 134 12 6633.75   # L134
Set-Cell $ws 134 13 -3755.4786   # M134
Set-Cell $ws 134 14 -11703.75   # N134
Set-Cell $ws 136 8 1325.9166   # H136
Set-Cell $ws 136 9 1341.4736   # I136
Set-Cell $ws 136 10 1266.8   # J136
Set-Cell $ws 136 11 4024.4208   # K136
Set-Cell $ws 136 12 3800.4   # L136
Set-Cell $ws 136 13 -1474.4208   # M136
Set-Cell $ws 136 14 -8900.4   # N136

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
Set-Cell $ws 5 8 1710.0667   # H5
Set-Cell $ws 5 9 1589.6   # I5
Set-Cell $ws 5 10 1951   # J5
Set-Cell $ws 5 11 4768.799999999999   # K5
Set-Cell $ws 5 12 5853   # L5
Set-Cell $ws 5 13 -4656.799999999999   # M5
Set-Cell $ws 5 14 -6077   # N5
Set-Cell $ws 44 8 3491   # H44
Set-Cell $ws 44 9 3000   # I44
Set-Cell $ws 44 10 3561.1428   # J44
Set-Cell $ws 44 11 9000   # K44
Set-Cell $ws 44 12 10683.4284   # L44
Set-Cell $ws 44 13 -8602   # M44
Set-Cell $ws 44 14 -11479.4284   # N44
Set-Cell $ws 56 8 7764   # H56
Set-Cell $ws 56 9 7764   # I56
Set-Cell $ws 56 10 0   # J56
Set-Cell $ws 56 11 7764   # K56
Set-Cell $ws 56 12 0   # L56
Set-Cell $ws 56 13 -7234   # M56
Set-Cell $ws 107 8 2140.45   # H107
Set-Cell $ws 107 9 6149.8   # I107
Set-Cell $ws 107 10 804   # J107
Set-Cell $ws 107 11 18449.4   # K107
Set-Cell $ws 107 12 2412   # L107
Set-Cell $ws 107 13 -16529.4   # M107
Set-Cell $ws 107 14 -6252   # N107
Set-Cell $ws 129 8 1936.4706   # H129
Set-Cell $ws 129 9 1335.8334   # I129
Set-Cell $ws 129 10 3378   # J129
Set-Cell $ws 129 11 4007.5002   # K129
Set-Cell $ws 129 12 10134   # L129
Set-Cell $ws 129 13 992.4998000000001   # M129
Set-Cell $ws 129 14 -20134   # N129
Set-Cell $ws 131 8 134713.47   # H131
Set-Cell $ws 131 9 851603.6   # I131
Set-Cell $ws 131 10 1956.037   # J131
Set-Cell $ws 131 11 2554810.8   # K131
Set-Cell $ws 131 12 5868.111   # L131
Set-Cell $ws 131 13 -2549770.8   # M131
Set-Cell $ws 131 14 -15948.111   # N131
Set-Cell $ws 135 8 1710.0667   # H135
Set-Cell $ws 135 9 1589.6   # I135
Set-Cell $ws 135 10 1951   # J135
Set-Cell $ws 135 11 14306.4   # K135
Set-Cell $ws 135 12 17559   # L135
Set-Cell $ws 135 13 -11771.4   # M135
Set-Cell $ws 135 14 -22629   # N135
Set-Cell $ws 137 8 3799.4443   # H137
Set-Cell $ws 137 9 3799.4443   # I137
Set-Cell $ws 137 10 0   # J137
Set-Cell $ws 137 11 11398.3329   # K137
Set-Cell $ws 137 12 0   # L137
Set-Cell $ws 137 13 -6298.332900000001   # M137

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
Set-Cell $ws 70 8 7056.696   # H70
Set-Cell $ws 70 9 7926.8125   # I70
Set-Cell $ws 70 10 5067.857   # J70
Set-Cell $ws 70 11 7926.8125   # K70
Set-Cell $ws 70 12 5067.857   # L70
Set-Cell $ws 70 13 -7656.8125   # M70
Set-Cell $ws 70 14 -5607.857   # N70
Set-Cell $ws 73 8 7056.696   # H73
Set-Cell $ws 73 9 7926.8125   # I73
Set-Cell $ws 73 10 5067.857   # J73
Set-Cell $ws 73 11 7926.8125   # K73
Set-Cell $ws 73 12 5067.857   # L73
Set-Cell $ws 73 13 -6990.8125   # M73
Set-Cell $ws 73 14 -6939.857   # N73
Set-Cell $ws 80 8 18624.75   # H80
Set-Cell $ws 80 9 6999.8   # I80
Set-Cell $ws 80 10 37999.668   # J80
Set-Cell $ws 80 11 6999.8   # K80
Set-Cell $ws 80 12 37999.668   # L80
Set-Cell $ws 80 13 -6001.8   # M80
Set-Cell $ws 80 14 -39995.668   # N80
Set-Cell $ws 83 8 18624.75   # H83
Set-Cell $ws 83 9 6999.8   # I83
Set-Cell $ws 83 10 37999.668   # J83
Set-Cell $ws 83 11 34999   # K83
Set-Cell $ws 83 12 189998.34   # L83
Set-Cell $ws 83 13 -30007   # M83
Set-Cell $ws 83 14 -199982.34   # N83

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
Set-Cell $ws 55 8 1123.6666   # H55
Set-Cell $ws 55 9 87.166664   # I55
Set-Cell $ws 55 10 2505.6667   # J55
Set-Cell $ws 55 11 87.166664   # K55
Set-Cell $ws 55 12 2505.6667   # L55
Set-Cell $ws 55 13 85.833336   # M55
Set-Cell $ws 55 14 -2851.6667   # N55
Set-Cell $ws 82 8 4365   # H82
Set-Cell $ws 82 9 2718.625   # I82
Set-Cell $ws 82 10 6011.375   # J82
Set-Cell $ws 82 11 2718.625   # K82
Set-Cell $ws 82 12 6011.375   # L82
Set-Cell $ws 82 13 -2357.625   # M82
Set-Cell $ws 82 14 -6733.375   # N82
Set-Cell $ws 85 8 4365   # H85
Set-Cell $ws 85 9 2718.625   # I85
Set-Cell $ws 85 10 6011.375   # J85
Set-Cell $ws 85 11 2718.625   # K85
Set-Cell $ws 85 12 6011.375   # L85
Set-Cell $ws 85 13 -1470.625   # M85
Set-Cell $ws 85 14 -8507.375   # N85
Set-Cell $ws 104 8 7947.8335   # H104
Set-Cell $ws 104 9 0   # I104
Set-Cell $ws 104 10 7947.8335   # J104
Set-Cell $ws 104 11 0   # K104
Set-Cell $ws 104 12 7947.8335   # L104
Set-Cell $ws 104 14 -14935.8335   # N104
